$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lgi3", "Adam23", "ECs", 1, 0.3333333333333333, 0.04481866666666667, 0.134456, 0.02406616871877758, 0.02406616871877757, 2, 0.6666666666666666, 0.2320676666666667, 0.696203, 0.01170834749781651, 0.01170834749781651, 0.01040096339644444, 0.093608670568, 0.0002817750663005295, 0.0002817750663005294),
    @("ECs", "Lgi3", "Adam23", "FAPs", 1, 0.3333333333333333, 0.04481866666666667, 0.134456, 0.02406616871877758, 0.02406616871877757, 3, 1, 17.63262733333333, 52.897882, 0.889606600882922, 0.8896066008829221, 0.7902708469102222, 7.112437622191999, 0.02140942255018663, 0.02140942255018663),
    @("ECs", "Lgi3", "Adam23", "Inflammatory-Mac", 1, 0.3333333333333333, 0.04481866666666667, 0.134456, 0.02406616871877758, 0.02406616871877757, 3, 1, 0.7378426666666668, 2.213528, 0.03722585944063267, 0.03722585944063268, 0.03306912452977778, 0.297622120768, 0.0008958838139997649, 0.0008958838139997649),
    @("ECs", "Lgi3", "Adam23", "MuSCs", 1, 0.3333333333333333, 0.04481866666666667, 0.134456, 0.02406616871877758, 0.02406616871877757, 3, 1, 1.143822333333333, 3.431467, 0.05770846730521116, 0.05770846730521117, 0.05126459188355555, 0.461381326952, 0.001388821710669271, 0.001388821710669271),
    @("ECs", "Lgi3", "Adam23", "Resolving-Mac", 1, 0.3333333333333333, 0.04481866666666667, 0.134456, 0.02406616871877758, 0.02406616871877757, 3, 1, 0.074342, 0.223026, 0.003750724873417703, 0.003750724873417703, 0.003331909317333334, 0.029987183856, 0.00009026557762138612, 0.00009026557762138611),
    @("FAPs", "Lgi3", "Adam23", "ECs", 3, 1, 1.606609333333333, 4.819827999999999, 0.8626970447097064, 0.8626970447097063, 2, 0.6666666666666666, 0.2320676666666667, 0.696203, 0.01170834749781651, 0.01170834749781651, 0.3728420792315555, 3.355578713084, 0.01010075678480059, 0.01010075678480059),
    @("FAPs", "Lgi3", "Adam23", "FAPs", 3, 1, 1.606609333333333, 4.819827999999999, 0.8626970447097064, 0.8626970447097063, 3, 1, 17.63262733333333, 52.897882, 0.889606600882922, 0.8896066008829221, 28.32874364492177, 254.9586928042959, 0.7674609855359441, 0.7674609855359441),
    @("FAPs", "Lgi3", "Adam23", "Inflammatory-Mac", 3, 1, 1.606609333333333, 4.819827999999999, 0.8626970447097064, 0.8626970447097063, 3, 1, 0.7378426666666668, 2.213528, 0.03722585944063267, 0.03722585944063268, 1.185424914798222, 10.668824233184, 0.03211463892621273, 0.03211463892621273),
    @("FAPs", "Lgi3", "Adam23", "MuSCs", 3, 1, 1.606609333333333, 4.819827999999999, 0.8626970447097064, 0.8626970447097063, 3, 1, 1.143822333333333, 3.431467, 0.05770846730521116, 0.05770846730521117, 1.837675636408444, 16.539080727676, 0.04978492419893239, 0.04978492419893239),
    @("FAPs", "Lgi3", "Adam23", "Resolving-Mac", 3, 1, 1.606609333333333, 4.819827999999999, 0.8626970447097064, 0.8626970447097063, 3, 1, 0.074342, 0.223026, 0.003750724873417703, 0.003750724873417703, 0.1194385510586667, 1.074946959528, 0.00323573926381664, 0.00323573926381664),
    @("MuSCs", "Lgi3", "Adam23", "ECs", 3, 1, 0.210882, 0.6326459999999999, 0.113236786571516, 0.113236786571516, 2, 0.6666666666666666, 0.2320676666666667, 0.696203, 0.01170834749781651, 0.01170834749781651, 0.048938893682, 0.440450043138, 0.001325815646715392, 0.001325815646715392),
    @("MuSCs", "Lgi3", "Adam23", "FAPs", 3, 1, 0.210882, 0.6326459999999999, 0.113236786571516, 0.113236786571516, 3, 1, 17.63262733333333, 52.897882, 0.889606600882922, 0.8896066008829221, 3.718403717308, 33.46563345577199, 0.1007361927967913, 0.1007361927967913),
    @("MuSCs", "Lgi3", "Adam23", "Inflammatory-Mac", 3, 1, 0.210882, 0.6326459999999999, 0.113236786571516, 0.113236786571516, 3, 1, 0.7378426666666668, 2.213528, 0.03722585944063267, 0.03722585944063268, 0.155597737232, 1.400379635088, 0.004215336700420177, 0.004215336700420177),
    @("MuSCs", "Lgi3", "Adam23", "MuSCs", 3, 1, 0.210882, 0.6326459999999999, 0.113236786571516, 0.113236786571516, 3, 1, 1.143822333333333, 3.431467, 0.05770846730521116, 0.05770846730521117, 0.241211541298, 2.170903871682, 0.006534721395609507, 0.006534721395609506),
    @("MuSCs", "Lgi3", "Adam23", "Resolving-Mac", 3, 1, 0.210882, 0.6326459999999999, 0.113236786571516, 0.113236786571516, 3, 1, 0.074342, 0.223026, 0.003750724873417703, 0.003750724873417703, 0.015677389644, 0.141096506796, 0.0004247200319796769, 0.0004247200319796768)
)

$r = 2
foreach ($rowvals in $data) {
    $c = 1
    foreach ($v in $rowvals) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
    $r = $r + 1
}

Write-Host "done"